$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the Title text in row 2 (remove the "CRs- " prefix)
$ws.Range("B2").Value = "CR list filter bar in the top menu in Admin mode"

# 2. Update the "Updated On" timestamp for row 2
$ws.Range("AD2").Value = "9/28/2023 1:35 PM"

# 3. Add a new "Forecast" value for row 5 (previously empty)
$ws.Range("H5").Value = "530s"
